$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the formatting of the
# existing header cells (bold font, thin border, centered alignment)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data values for rows 2 and 3
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 5
